$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 0.08902255639097745
$ws.Range("G2").Value = "JP"

$ws.Range("F3").Value = "kraemer (1980)"
$ws.Range("F4").Value = "kraemer (1980)"
$ws.Range("F5").Value = "kraemer (1980)"
$ws.Range("F6").Value = "kraemer (1980)"

$ws.Range("B7").Value = 0.5837902834806478
$ws.Range("C7").Value = 0.94862625852474
$ws.Range("D7").Value = 0.7850223480679014
$ws.Range("F7").Value = "kraemer (1980)"

$ws.Range("F8").Value = "kraemer (1980)"
$ws.Range("F9").Value = "kraemer (1980)"
$ws.Range("F10").Value = "kraemer (1980)"
$ws.Range("F11").Value = "kraemer (1980)"
$ws.Range("F12").Value = "kraemer (1980)"
$ws.Range("F13").Value = "kraemer (1980)"
$ws.Range("F14").Value = "kraemer (1980)"

$ws.Range("E16").ClearContents()
$ws.Range("F16").Value = "kraemer (1980)"
$ws.Range("G16").ClearContents()

$ws.Range("F17").Value = "kraemer (1980)"

$ws.Range("F19").Value = "kraemer (1980)"
$ws.Range("F20").Value = "kraemer (1980)"
$ws.Range("F21").Value = "kraemer (1980)"

$ws.Range("F34").Value = "kraemer (1980)"
$ws.Range("F35").Value = "kraemer (1980)"
$ws.Range("F36").Value = "kraemer (1980)"
